$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B5: ORGANIZATION/PERSON value trimmed down ---
$ws.Range("B5").Value = "Dugal Harris"

# --- Row 8 (3.1a): clear the budget figures / notes that were removed ---
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("H8").ClearContents()

# --- Row 9 (3.1b): activities / deliverable / budget notes rewritten ---
$ws.Range("H9").Value = "3 months consultant time at ZAR 300/hr for 12 days work"
$ws.Range("B9").Value = "1) Apply differential correction to field DGPS data acquired during Q1&2.`n2) Process Q1&2 field allometric data to produce woody carbon stock (CS) estimates.`n3) Perform preliminary correlation analysis on data acquired to date. `n"
$ws.Range("C9").Value = "1) Results of activ ities 1-3 (models, graphs and correlation statistics).  "
$ws.Rows.Item(9).RowHeight = 93.75

# --- Row 10 (3.1d): clear the notes/budget figures that were removed ---
$ws.Range("B10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("H10").ClearContents()

# --- Update the active cell selection on the frozen pane ---
$ws.Range("B17").Select()
